$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.932.99'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '2.561.33'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''301.98'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').Value = '''92.54'
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').Value = '''36.14'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '''7.71'
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('E13').Value = '  +6.88%  '
$ws.Range('D14').Value = '2.549.61'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '42.982.44'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '0.0₃0995'
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('D19').Value = '''12.69'
$ws.Range('E19').Value = '  +3.12%  '
$ws.Range('D20').Value = '''6.58'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').Value = '''71.63'
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('D22').Value = '''253.34'
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').Value = '''2.94'
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('E24').Value = '  -3.88%  '
$ws.Range('E25').Value = '  -2.01%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').Value = '''10.29'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('D28').Value = '''37.17'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('E29').Value = '  -3.76%  '
$ws.Range('D30').Value = '''6.03'
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('D31').Value = '''153.85'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  -1.54%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = '''3.38'
$ws.Range('E33').Value = '  -5.82%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''2.14'
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('D35').Value = '''0.0800'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').Value = '''17.99'
$ws.Range('E36').Value = '  +7.07%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '''0.114'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').Value = '''23.07'
$ws.Range('E39').Value = '  -5.48%  '
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('D41').Value = '''0.0311'
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('D43').Value = '2.098.15'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('D44').Value = '''2.01'
$ws.Range('E44').Value = '  +25.72%  '
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').Value = '''9.28'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').Value = '''85.12'
$ws.Range('E47').Value = '  -3.11%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '''75.98'
$ws.Range('E48').Value = '  +10.63%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''106.65'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('D50').Value = '2.811.39'
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('E51').Value = '  +2.05%  '
